$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Table Entry")

$ws.Range("A4").Value  = "ui9v547jxpt9of6"
$ws.Range("A5").Value  = "6xmay00ev2l4pl1"
$ws.Range("A6").Value  = "1f1s9j33q3lll85"
$ws.Range("A7").Value  = "rdo4e4x959n4262"
$ws.Range("A9").Value  = "ka5njrmwml2tw7y"
$ws.Range("A10").Value = "225el2d7iyhk1i6"
$ws.Range("A11").Value = "4f1z9ba31h53526"
$ws.Range("A13").Value = "ds2d4310f2yqd8a"
